# Regenerate the "K" (strikeouts) column (column G) values in the save_data
# sheet for velasquez_vince.xlsx. The values were recomputed upstream (switch
# from "Strike#" to "K") and need to be written back into the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$kValues = @{
    2  = 5
    3  = 1
    4  = 5
    5  = 2
    6  = 4
    7  = 1
    8  = 2
    9  = 1
    10 = 0
    11 = 2
    12 = 1
    13 = 2
    14 = 1
    15 = 0
    16 = 1
    17 = 4
    18 = 0
    19 = 1
    20 = 3
    21 = 5
    22 = 3
    23 = 4
    24 = 2
    25 = 6
    26 = 6
    27 = 5
    28 = 2
    29 = 5
    31 = 2
    32 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
